$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.403.14"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.572.62"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3766"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.83"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3419"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07637"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.148"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.011"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.942"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "1.574.21"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001135"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.216"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "22.410.92"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.666"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -10.88%  "
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.036"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "1.750.83"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.138"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.009"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9849"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.86%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08511"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02541"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.368"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2316"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06524"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.443"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6375"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.793"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5971"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.096"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.283"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07327"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.39%  "
